$wb = $excel.ActiveWorkbook

# --- 1. Update the Q1 ("ms_source") comment text on the "Export as TSV" sheet ---
$wsMain = $wb.Worksheets.Item("Export as TSV")
$comment = $wsMain.Range("Q1").Comment
$comment.Text("The ion source type used for surface sampling (MALDI, MALDI-2, DESI, nanoDESI or SIMS).")

# --- 2. Remove "ESI" from the "ms_source list" sheet (A7, shifting nanoDESI up) ---
$wsSource = $wb.Worksheets.Item("ms_source list")
$wsSource.Range("A7").EntireRow.Delete()

# --- 3. Fix the ms_source data validation range on column Q (was A1:A8, now A1:A7) ---
$dv = $wsMain.Range("Q2:Q1048576").Validation
$dv.Modify(3, 1, 1, "'ms_source list'!`$A`$1:`$A`$7")
